$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 44306
$ws.Range("M2").Value = 45
$ws.Range("D3").Value = 44333
$ws.Range("L3").Value = "Especial"
$ws.Range("M3").Value = 58
$ws.Range("N3").Value = 10000
$ws.Range("O3").Value = 10000
$ws.Range("P3").Value = 10000
$ws.Range("R3").Value = "Provincia de Quillota"
$ws.Range("S3").Value = 1000
$ws.Range("D4").Value = 44333
$ws.Range("L4").Value = "Primera"
$ws.Range("M4").Value = 65
$ws.Range("N4").Value = 9000
$ws.Range("O4").Value = 9000
$ws.Range("P4").Value = 9000
$ws.Range("R4").Value = "Provincia de Quillota"
$ws.Range("S4").Value = 900
$ws.Range("D5").Value = 44333
$ws.Range("L5").Value = "Segunda"
$ws.Range("M5").Value = 60
$ws.Range("N5").Value = 8000
$ws.Range("O5").Value = 8000
$ws.Range("P5").Value = 8000
$ws.Range("S5").Value = 800
$ws.Range("D6").Value = 44321
$ws.Range("M6").Value = 58
$ws.Range("N6").Value = 9000
$ws.Range("O6").Value = 9000
$ws.Range("P6").Value = 9000
$ws.Range("S6").Value = 900
$ws.Range("D7").Value = 44307
$ws.Range("L7").Value = "Primera"
$ws.Range("M7").Value = 40
$ws.Range("N7").Value = 10000
$ws.Range("O7").Value = 10000
$ws.Range("P7").Value = 10000
$ws.Range("S7").Value = 1000
$ws.Range("D8").Value = 44328
$ws.Range("M8").Value = 45
$ws.Range("N8").Value = 8000
$ws.Range("O8").Value = 8000
$ws.Range("P8").Value = 8000
$ws.Range("S8").Value = 800
$ws.Range("D9").Value = 44328
$ws.Range("L9").Value = "Segunda"
$ws.Range("M9").Value = 48
$ws.Range("N9").Value = 7000
$ws.Range("O9").Value = 7000
$ws.Range("P9").Value = 7000
$ws.Range("S9").Value = 700
$ws.Range("D10").Value = 44319
$ws.Range("L10").Value = "Primera"
$ws.Range("M10").Value = 68
$ws.Range("D11").Value = 44319
$ws.Range("L11").Value = "Segunda"
$ws.Range("M11").Value = 57
$ws.Range("N11").Value = 8000
$ws.Range("O11").Value = 8000
$ws.Range("P11").Value = 8000
$ws.Range("S11").Value = 800
$ws.Range("D12").Value = 44312
$ws.Range("L12").Value = "Primera"
$ws.Range("M12").Value = 48
$ws.Range("N12").Value = 10000
$ws.Range("O12").Value = 10000
$ws.Range("P12").Value = 10000
$ws.Range("S12").Value = 1000
$ws.Range("D13").Value = 44699
$ws.Range("L13").Value = "Especial"
$ws.Range("M13").Value = 56
$ws.Range("N13").Value = 12000
$ws.Range("O13").Value = 12000
$ws.Range("P13").Value = 12000
$ws.Range("S13").Value = 1200
$ws.Range("D14").Value = 44699
$ws.Range("M14").Value = 60
$ws.Range("D15").Value = 44309
$ws.Range("D16").Value = 44314
$ws.Range("M16").Value = 47
$ws.Range("N16").Value = 9000
$ws.Range("O16").Value = 9000
$ws.Range("P16").Value = 9000
$ws.Range("S16").Value = 900
$ws.Range("D17").Value = 44343
$ws.Range("L17").Value = "Especial"
$ws.Range("M17").Value = 47
$ws.Range("N17").Value = 10000
$ws.Range("O17").Value = 10000
$ws.Range("P17").Value = 10000
$ws.Range("R17").Value = "Región Metropolitana"
$ws.Range("S17").Value = 1000
$ws.Range("D18").Value = 44343
$ws.Range("M18").Value = 50
$ws.Range("N18").Value = 9000
$ws.Range("O18").Value = 9000
$ws.Range("P18").Value = 9000
$ws.Range("R18").Value = "Región Metropolitana"
$ws.Range("S18").Value = 900
$ws.Range("D19").Value = 44343
$ws.Range("M19").Value = 58
$ws.Range("R19").Value = "Región Metropolitana"
$ws.Range("D20").Value = 44326
$ws.Range("M20").Value = 65
$ws.Range("D21").Value = 44326
$ws.Range("M21").Value = 67
$ws.Range("D22").Value = 44308
$ws.Range("D23").Value = 44308
$ws.Range("L23").Value = "Segunda"
$ws.Range("M23").Value = 48
$ws.Range("N23").Value = 8000
$ws.Range("O23").Value = 8000
$ws.Range("P23").Value = 8000
$ws.Range("S23").Value = 800
$ws.Range("D24").Value = 44323
$ws.Range("D25").Value = 44323
$ws.Range("L25").Value = "Segunda"
$ws.Range("M25").Value = 50
$ws.Range("N25").Value = 9000
$ws.Range("O25").Value = 9000
$ws.Range("P25").Value = 9000
$ws.Range("S25").Value = 900
$ws.Range("D26").Value = 44315
$ws.Range("M26").Value = 45
$ws.Range("D27").Value = 44329
$ws.Range("L27").Value = "Primera"
$ws.Range("M27").Value = 56
$ws.Range("N27").Value = 9000
$ws.Range("O27").Value = 9000
$ws.Range("P27").Value = 9000
$ws.Range("R27").Value = "Región Metropolitana"
$ws.Range("S27").Value = 900
$ws.Range("D28").Value = 44329
$ws.Range("L28").Value = "Segunda"
$ws.Range("M28").Value = 50
$ws.Range("N28").Value = 8000
$ws.Range("O28").Value = 8000
$ws.Range("P28").Value = 8000
$ws.Range("S28").Value = 800
$ws.Range("D29").Value = 44322
$ws.Range("M29").Value = 56
$ws.Range("N29").Value = 10000
$ws.Range("O29").Value = 10000
$ws.Range("P29").Value = 10000
$ws.Range("R29").Value = "Provincia de Quillota"
$ws.Range("S29").Value = 1000
$ws.Range("D30").Value = 44322
$ws.Range("M30").Value = 40
$ws.Range("R30").Value = "Provincia de Quillota"
$ws.Range("D31").Value = 44301
$ws.Range("N31").Value = 10000
$ws.Range("O31").Value = 10000
$ws.Range("P31").Value = 10000
$ws.Range("S31").Value = 1000
$ws.Range("D32").Value = 44302
$ws.Range("L32").Value = "Primera"
$ws.Range("M32").Value = 45
$ws.Range("N32").Value = 10000
$ws.Range("O32").Value = 10000
$ws.Range("P32").Value = 10000
$ws.Range("S32").Value = 1000
